$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2024-09-03 12:18:17"
$ws.Range("B4").Value = "ORM-0515835"
$ws.Range("C4").Value = "Z195631098"
$ws.Range("D4").Value = "C18H19CL2N5OS"
$ws.Range("E4").Value = "Duplicate"
$ws.Range("F4").Value = "/home/robekott/ERAT/examples/compound_test.sdf"

$ws.Range("A5").Value = "2024-09-03 12:18:18"
$ws.Range("B5").Value = "ORM-0515836"
$ws.Range("C5").Value = "Z2754556176"
$ws.Range("D5").Value = "C17H28N4O2"
$ws.Range("E5").Value = "Duplicate"
$ws.Range("F5").Value = "/home/robekott/ERAT/examples/compound_test.sdf"
